# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 07:52"

# Pakistan (row 15)
$ws.Range("B15").Value = 248872
$ws.Range("C15").Value = 2521
$ws.Range("D15").Value = 156700
$ws.Range("E15").Value = 86975
$ws.Range("G15").Value = 74
$ws.Range("H15").Value = 5197

# Belgica (row 33)
$ws.Range("B33").Value = 62606
$ws.Range("C33").Value = 137
$ws.Range("E33").Value = 35628

# Kirguistan (row 72)
$ws.Range("B72").Value = 10629
$ws.Range("C72").Value = 719
$ws.Range("D72").Value = 3387
$ws.Range("E72").Value = 7110
$ws.Range("G72").Value = 7
$ws.Range("H72").Value = 132

# Tailandia (row 103)
$ws.Range("B103").Value = 3217
$ws.Range("C103").Value = 1
$ws.Range("E103").Value = 71

# Butan (row 187)
$ws.Range("D187").Value = 76
$ws.Range("E187").Value = 6
